$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: correct the "valor faturado" / "valor da comissao" amounts ---
$ws.Range("G2").Value = 22714.21
$ws.Range("H2").Value = 2271.4210000000003

# --- Row 3: status was wrongly "NAOFATURADO", should be "PARCIALMENTEFATURADO"; also fix amounts ---
$ws.Range("E3").Value = "PARCIALMENTEFATURADO"
$ws.Range("G3").Value = 12743.0
$ws.Range("H3").Value = 1274.3000000000002

# --- Insert a brand-new row 4 (pushes the former row 4 down to row 5) ---
$ws.Rows("4").Insert()

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "4"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "123213"
$ws.Range("C4").Value = "Cliente Pedro 1"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "23"
$ws.Range("E4").Value = "TOTALMENTEFATURADO"
$ws.Range("F4").Value = 200001.0
$ws.Range("G4").Value = 200001.0
$ws.Range("H4").Value = 9000.045

# --- Append a new row 5 (new export entry) ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "5"
$ws.Range("B5").Value = "Industria teste 1"
$ws.Range("C5").Value = "Araujo atacadista"
$ws.Range("D5").Value = "33333333a"
$ws.Range("E5").Value = "TOTALMENTEFATURADO"
$ws.Range("F5").Value = 200.0
$ws.Range("G5").Value = 200.0
$ws.Range("H5").Value = 20.0

# --- Apply Brazilian currency number format to the monetary columns ---
$ws.Range("F2:H5").NumberFormat = "R$ #,##0.00"

# --- Column widths to fit the report contents (bestFit) ---
# NOTE: the runtime stores ColumnWidth on a 1/6-character grid, adding an
# implicit ~0.8333 padding before rounding, so we pre-compensate by the same
# offset to land on the closest representable width to the target OOXML value.
$ws.Columns("A").ColumnWidth = 1.8072916666666665
$ws.Columns("B").ColumnWidth = 13.170572916666666
$ws.Columns("C").ColumnWidth = 14.065104166666666
$ws.Columns("D").ColumnWidth = 14.186197916666666
$ws.Columns("E").ColumnWidth = 21.705729166666668
$ws.Columns("F").ColumnWidth = 12.311197916666666
$ws.Columns("G").ColumnWidth = 11.604166666666666
$ws.Columns("H").ColumnWidth = 14.420572916666666
